# Auto-generated edit script: updates market-price / profit columns
# (H..N) for specific Leve rows across all 8 job sheets, matching the
# scheduled-runner data refresh described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Cells.Item(15, 8).Value = 285.02
$ws.Cells.Item(15, 9).Value = 285.02
$ws.Cells.Item(15, 11).Value = 855.0599999999999
$ws.Cells.Item(15, 13).Value = -686.0599999999999
# ALC row 18
$ws.Cells.Item(18, 8).Value = 928.13336
$ws.Cells.Item(18, 9).Value = 980
$ws.Cells.Item(18, 10).Value = 202
$ws.Cells.Item(18, 11).Value = 980
$ws.Cells.Item(18, 12).Value = 202
$ws.Cells.Item(18, 13).Value = -696
$ws.Cells.Item(18, 14).Value = -770
# ALC row 19
$ws.Cells.Item(19, 8).Value = 4311545.5
$ws.Cells.Item(19, 9).Value = 6944885.5
$ws.Cells.Item(19, 10).Value = 2443.818
$ws.Cells.Item(19, 11).Value = 6944885.5
$ws.Cells.Item(19, 12).Value = 2443.818
$ws.Cells.Item(19, 13).Value = -6944710.5
$ws.Cells.Item(19, 14).Value = -2793.818
# ALC row 62
$ws.Cells.Item(62, 8).Value = 76090290
$ws.Cells.Item(62, 9).Value = 25003830
$ws.Cells.Item(62, 10).Value = 416666660
$ws.Cells.Item(62, 11).Value = 25003830
$ws.Cells.Item(62, 12).Value = 416666660
$ws.Cells.Item(62, 13).Value = -25003206
$ws.Cells.Item(62, 14).Value = -416667908
# ALC row 65
$ws.Cells.Item(65, 8).Value = 76090290
$ws.Cells.Item(65, 9).Value = 25003830
$ws.Cells.Item(65, 10).Value = 416666660
$ws.Cells.Item(65, 11).Value = 125019150
$ws.Cells.Item(65, 12).Value = 2083333300
$ws.Cells.Item(65, 13).Value = -125016030
$ws.Cells.Item(65, 14).Value = -2083339540
# ALC row 113
$ws.Cells.Item(113, 8).Value = 4002743.2
$ws.Cells.Item(113, 9).Value = 8335273.5
$ws.Cells.Item(113, 10).Value = 3484.6155
$ws.Cells.Item(113, 11).Value = 8335273.5
$ws.Cells.Item(113, 12).Value = 3484.6155
$ws.Cells.Item(113, 13).Value = -8332019.5
$ws.Cells.Item(113, 14).Value = -9992.6155
# ALC row 118
$ws.Cells.Item(118, 8).Value = 2201.9092
$ws.Cells.Item(118, 10).Value = 3984
$ws.Cells.Item(118, 12).Value = 11952
$ws.Cells.Item(118, 14).Value = -15266
# ALC row 127
$ws.Cells.Item(127, 8).Value = 1258.0416
$ws.Cells.Item(127, 9).Value = 452.52942
$ws.Cells.Item(127, 10).Value = 3214.2856
$ws.Cells.Item(127, 11).Value = 1357.58826
$ws.Cells.Item(127, 12).Value = 9642.856800000001
$ws.Cells.Item(127, 13).Value = 3602.41174
$ws.Cells.Item(127, 14).Value = -19562.8568
# ALC row 132
$ws.Cells.Item(132, 8).Value = 3493289
$ws.Cells.Item(132, 9).Value = 717833.2
$ws.Cells.Item(132, 10).Value = 27778528
$ws.Cells.Item(132, 11).Value = 2153499.6
$ws.Cells.Item(132, 12).Value = 83335584
$ws.Cells.Item(132, 13).Value = -2150969.6
$ws.Cells.Item(132, 14).Value = -83340644
# ALC row 137
$ws.Cells.Item(137, 8).Value = 27833710
$ws.Cells.Item(137, 9).Value = 62501500
$ws.Cells.Item(137, 10).Value = 17928628
$ws.Cells.Item(137, 11).Value = 187504500
$ws.Cells.Item(137, 12).Value = 53785884
$ws.Cells.Item(137, 13).Value = -187501950
$ws.Cells.Item(137, 14).Value = -53790984

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Cells.Item(2, 8).Value = 22862.5
$ws.Cells.Item(2, 9).Value = 25971.428
$ws.Cells.Item(2, 10).Value = 1100
$ws.Cells.Item(2, 11).Value = 25971.428
$ws.Cells.Item(2, 12).Value = 1100
$ws.Cells.Item(2, 13).Value = -25858.428
$ws.Cells.Item(2, 14).Value = -1326
# ARM row 45
$ws.Cells.Item(45, 8).Value = 2260.68
$ws.Cells.Item(45, 9).Value = 1281.4375
$ws.Cells.Item(45, 10).Value = 4001.5557
$ws.Cells.Item(45, 11).Value = 1281.4375
$ws.Cells.Item(45, 12).Value = 4001.5557
$ws.Cells.Item(45, 13).Value = -904.4375
$ws.Cells.Item(45, 14).Value = -4755.5557
# ARM row 74
$ws.Cells.Item(74, 8).Value = 45979330
$ws.Cells.Item(74, 9).Value = 47620268
$ws.Cells.Item(74, 10).Value = 41671852
$ws.Cells.Item(74, 11).Value = 47620268
$ws.Cells.Item(74, 12).Value = 41671852
$ws.Cells.Item(74, 13).Value = -47619394
$ws.Cells.Item(74, 14).Value = -41673600
# ARM row 77
$ws.Cells.Item(77, 8).Value = 45979330
$ws.Cells.Item(77, 9).Value = 47620268
$ws.Cells.Item(77, 10).Value = 41671852
$ws.Cells.Item(77, 11).Value = 238101340
$ws.Cells.Item(77, 12).Value = 208359260
$ws.Cells.Item(77, 13).Value = -238096972
$ws.Cells.Item(77, 14).Value = -208367996
# ARM row 116
$ws.Cells.Item(116, 8).Value = 22862.5
$ws.Cells.Item(116, 9).Value = 25971.428
$ws.Cells.Item(116, 10).Value = 1100
$ws.Cells.Item(116, 11).Value = 25971.428
$ws.Cells.Item(116, 12).Value = 1100
$ws.Cells.Item(116, 13).Value = -23677.428
$ws.Cells.Item(116, 14).Value = -5688
# ARM row 122
$ws.Cells.Item(122, 8).Value = 1654.8889
$ws.Cells.Item(122, 9).Value = 1726.7858
$ws.Cells.Item(122, 10).Value = 1403.25
$ws.Cells.Item(122, 11).Value = 5180.357400000001
$ws.Cells.Item(122, 12).Value = 4209.75
$ws.Cells.Item(122, 13).Value = -2730.357400000001
$ws.Cells.Item(122, 14).Value = -9109.75
# ARM row 135
$ws.Cells.Item(135, 8).Value = 45119.188
$ws.Cells.Item(135, 10).Value = 45119.188
$ws.Cells.Item(135, 12).Value = 45119.188
$ws.Cells.Item(135, 14).Value = -55259.188
# ARM row 138
$ws.Cells.Item(138, 8).Value = 45015.8
$ws.Cells.Item(138, 10).Value = 45015.8
$ws.Cells.Item(138, 12).Value = 45015.8
$ws.Cells.Item(138, 14).Value = -55295.8
# ARM row 139
$ws.Cells.Item(139, 8).Value = 30136.62
$ws.Cells.Item(139, 10).Value = 30136.62
$ws.Cells.Item(139, 12).Value = 30136.62
$ws.Cells.Item(139, 14).Value = -40416.62

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Cells.Item(3, 8).Value = 22862.5
$ws.Cells.Item(3, 9).Value = 25971.428
$ws.Cells.Item(3, 10).Value = 1100
$ws.Cells.Item(3, 11).Value = 25971.428
$ws.Cells.Item(3, 12).Value = 1100
$ws.Cells.Item(3, 13).Value = -25857.428
$ws.Cells.Item(3, 14).Value = -1328

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Cells.Item(16, 8).Value = 1950
$ws.Cells.Item(16, 9).Value = 1800
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 1800
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = -1513
$ws.Cells.Item(16, 14).Value = -2574
# CRP row 31
$ws.Cells.Item(31, 8).Value = 1464438.9
$ws.Cells.Item(31, 9).Value = 1852929.2
$ws.Cells.Item(31, 10).Value = 7600
$ws.Cells.Item(31, 11).Value = 1852929.2
$ws.Cells.Item(31, 12).Value = 7600
$ws.Cells.Item(31, 13).Value = -1852634.2
$ws.Cells.Item(31, 14).Value = -8190
# CRP row 34
$ws.Cells.Item(34, 8).Value = 1464438.9
$ws.Cells.Item(34, 9).Value = 1852929.2
$ws.Cells.Item(34, 10).Value = 7600
$ws.Cells.Item(34, 11).Value = 1852929.2
$ws.Cells.Item(34, 12).Value = 7600
$ws.Cells.Item(34, 13).Value = -1852727.2
$ws.Cells.Item(34, 14).Value = -8004
# CRP row 58
$ws.Cells.Item(58, 8).Value = 814283.3
$ws.Cells.Item(58, 9).Value = 3451.3057
$ws.Cells.Item(58, 10).Value = 2273780.8
$ws.Cells.Item(58, 11).Value = 3451.3057
$ws.Cells.Item(58, 12).Value = 2273780.8
$ws.Cells.Item(58, 13).Value = -3248.3057
$ws.Cells.Item(58, 14).Value = -2274186.8
# CRP row 100
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).ClearContents()
# CRP row 107
$ws.Cells.Item(107, 8).Value = 1067.3914
$ws.Cells.Item(107, 9).Value = 327.57144
$ws.Cells.Item(107, 10).Value = 1391.0625
$ws.Cells.Item(107, 11).Value = 327.57144
$ws.Cells.Item(107, 12).Value = 1391.0625
$ws.Cells.Item(107, 13).Value = 1592.42856
$ws.Cells.Item(107, 14).Value = -5231.0625
# CRP row 113
$ws.Cells.Item(113, 8).Value = 1950
$ws.Cells.Item(113, 9).Value = 1800
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 1800
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = 370
$ws.Cells.Item(113, 14).Value = -6340
# CRP row 122
$ws.Cells.Item(122, 8).Value = 5582.087
$ws.Cells.Item(122, 9).Value = 6189.4
$ws.Cells.Item(122, 10).Value = 1533.3334
$ws.Cells.Item(122, 11).Value = 18568.2
$ws.Cells.Item(122, 12).Value = 4600.0002
$ws.Cells.Item(122, 13).Value = -16118.2
$ws.Cells.Item(122, 14).Value = -9500.0002
# CRP row 132
$ws.Cells.Item(132, 8).Value = 1538.7812
$ws.Cells.Item(132, 9).Value = 1017.16
$ws.Cells.Item(132, 10).Value = 3401.7144
$ws.Cells.Item(132, 11).Value = 3051.48
$ws.Cells.Item(132, 12).Value = 10205.1432
$ws.Cells.Item(132, 13).Value = -521.48
$ws.Cells.Item(132, 14).Value = -15265.1432
# CRP row 136
$ws.Cells.Item(136, 8).Value = 814283.3
$ws.Cells.Item(136, 9).Value = 3451.3057
$ws.Cells.Item(136, 10).Value = 2273780.8
$ws.Cells.Item(136, 11).Value = 10353.9171
$ws.Cells.Item(136, 12).Value = 6821342.399999999
$ws.Cells.Item(136, 13).Value = -7803.917099999999
$ws.Cells.Item(136, 14).Value = -6826442.399999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 3
$ws.Cells.Item(3, 8).Value = 8711.444
$ws.Cells.Item(3, 9).Value = 7921.25
$ws.Cells.Item(3, 10).Value = 15033
$ws.Cells.Item(3, 11).Value = 23763.75
$ws.Cells.Item(3, 12).Value = 45099
$ws.Cells.Item(3, 13).Value = -23651.75
$ws.Cells.Item(3, 14).Value = -45323
# CUL row 8
$ws.Cells.Item(8, 8).Value = 249.8
$ws.Cells.Item(8, 9).Value = 249.8
$ws.Cells.Item(8, 11).Value = 749.4000000000001
$ws.Cells.Item(8, 13).Value = -610.4000000000001
# CUL row 121
$ws.Cells.Item(121, 8).Value = 2604906
$ws.Cells.Item(121, 9).Value = 413.2
$ws.Cells.Item(121, 10).Value = 3125804.5
$ws.Cells.Item(121, 11).Value = 1239.6
$ws.Cells.Item(121, 12).Value = 9377413.5
$ws.Cells.Item(121, 13).Value = 70.40000000000009
$ws.Cells.Item(121, 14).Value = -9380033.5
# CUL row 134
$ws.Cells.Item(134, 8).Value = 2576.2222
$ws.Cells.Item(134, 9).Value = 1588.5714
$ws.Cells.Item(134, 10).Value = 6033
$ws.Cells.Item(134, 11).Value = 4765.7142
$ws.Cells.Item(134, 12).Value = 18099
$ws.Cells.Item(134, 13).Value = 304.2857999999997
$ws.Cells.Item(134, 14).Value = -28239

$ws = $wb.Worksheets.Item("GSM")
# GSM row 101
$ws.Cells.Item(101, 8).Value = 35000
$ws.Cells.Item(101, 10).Value = 35000
$ws.Cells.Item(101, 12).Value = 35000
$ws.Cells.Item(101, 14).Value = -41490
# GSM row 126
$ws.Cells.Item(126, 8).Value = 8314.157999999999
$ws.Cells.Item(126, 9).Value = 18496.834
$ws.Cells.Item(126, 11).Value = 55490.50199999999
$ws.Cells.Item(126, 13).Value = -53020.50199999999
# GSM row 138
$ws.Cells.Item(138, 8).Value = 39345.105
$ws.Cells.Item(138, 10).Value = 39345.105
$ws.Cells.Item(138, 12).Value = 39345.105
$ws.Cells.Item(138, 14).Value = -49625.105

$ws = $wb.Worksheets.Item("LTW")
# LTW row 46
$ws.Cells.Item(46, 8).Value = 1030.5264
$ws.Cells.Item(46, 9).Value = 752.63635
$ws.Cells.Item(46, 10).Value = 1412.625
$ws.Cells.Item(46, 11).Value = 752.63635
$ws.Cells.Item(46, 12).Value = 1412.625
$ws.Cells.Item(46, 13).Value = -564.63635
$ws.Cells.Item(46, 14).Value = -1788.625
# LTW row 122
$ws.Cells.Item(122, 8).Value = 18853420
$ws.Cells.Item(122, 9).Value = 11797109
$ws.Cells.Item(122, 10).Value = 100001000
$ws.Cells.Item(122, 11).Value = 35391327
$ws.Cells.Item(122, 12).Value = 300003000
$ws.Cells.Item(122, 13).Value = -35388877
$ws.Cells.Item(122, 14).Value = -300007900
# LTW row 136
$ws.Cells.Item(136, 8).Value = 3573620.8
$ws.Cells.Item(136, 9).Value = 6252145.5
$ws.Cells.Item(136, 10).Value = 2254.6667
$ws.Cells.Item(136, 11).Value = 18756436.5
$ws.Cells.Item(136, 12).Value = 6764.000100000001
$ws.Cells.Item(136, 13).Value = -18753886.5
$ws.Cells.Item(136, 14).Value = -11864.0001

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Cells.Item(122, 8).Value = 5557462.5
$ws.Cells.Item(122, 9).Value = 13890166
$ws.Cells.Item(122, 10).Value = 2326.5833
$ws.Cells.Item(122, 11).Value = 41670498
$ws.Cells.Item(122, 12).Value = 6979.749899999999
$ws.Cells.Item(122, 13).Value = -41668048
$ws.Cells.Item(122, 14).Value = -11879.7499
# WVR row 136
$ws.Cells.Item(136, 8).Value = 1930.0938
$ws.Cells.Item(136, 9).Value = 1503.5555
$ws.Cells.Item(136, 10).Value = 2478.5
$ws.Cells.Item(136, 11).Value = 1930.0938
$ws.Cells.Item(136, 12).Value = 7435.5
$ws.Cells.Item(136, 13).Value = -1960.666499999999
$ws.Cells.Item(136, 14).Value = -12535.5

